# Auto-generated script: apply scheduled-runner market-data refresh
# to the Leve profit calculation sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    @{Cell="H15"; Value=2100.6377}
    @{Cell="I15"; Value=2100.6377}
    @{Cell="K15"; Value=6301.913100000001}
    @{Cell="M15"; Value=-6132.913100000001}
    @{Cell="H19"; Value=38462890}
    @{Cell="I19"; Value=1109.75}
    @{Cell="J19"; Value=55557012}
    @{Cell="K19"; Value=1109.75}
    @{Cell="L19"; Value=55557012}
    @{Cell="M19"; Value=-934.75}
    @{Cell="N19"; Value=-55557362}
    @{Cell="H96"; Value=1587.6666}
    @{Cell="I96"; Value=1587.6666}
    @{Cell="K96"; Value=4762.9998}
    @{Cell="M96"; Value=-3389.9998}
    @{Cell="H103"; Value=1508.75}
    @{Cell="I103"; Value=1275}
    @{Cell="K103"; Value=3825}
    @{Cell="M103"; Value=-3239}
    @{Cell="H129"; Value=887.5454999999999}
    @{Cell="I129"; Value=435.6}
    @{Cell="J129"; Value=1264.1666}
    @{Cell="K129"; Value=1306.8}
    @{Cell="L129"; Value=3792.4998}
    @{Cell="M129"; Value=3693.2}
    @{Cell="N129"; Value=-13792.4998}
)
foreach ($u in $updates) { $ws.Range($u.Cell).Value = $u.Value }

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    @{Cell="H45"; Value=1747.4736}
    @{Cell="I45"; Value=1747.6}
    @{Cell="J45"; Value=1747.2307}
    @{Cell="K45"; Value=1747.6}
    @{Cell="L45"; Value=1747.2307}
    @{Cell="M45"; Value=-1370.6}
    @{Cell="N45"; Value=-2501.2307}
    @{Cell="H122"; Value=3206217.8}
    @{Cell="I122"; Value=1412.875}
    @{Cell="J122"; Value=4033264.2}
    @{Cell="K122"; Value=4238.625}
    @{Cell="L122"; Value=12099792.6}
    @{Cell="M122"; Value=-1788.625}
    @{Cell="N122"; Value=-12104692.6}
    @{Cell="H132"; Value=2312.652}
    @{Cell="I132"; Value=1918.7297}
    @{Cell="K132"; Value=5756.189100000001}
    @{Cell="M132"; Value=-3226.189100000001}
)
foreach ($u in $updates) { $ws.Range($u.Cell).Value = $u.Value }

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    @{Cell="H99"; Value=1780.3684}
    @{Cell="I99"; Value=1077.5834}
    @{Cell="J99"; Value=2985.1428}
    @{Cell="K99"; Value=1077.5834}
    @{Cell="L99"; Value=2985.1428}
    @{Cell="M99"; Value=420.4166}
    @{Cell="N99"; Value=-5981.1428}
)
foreach ($u in $updates) { $ws.Range($u.Cell).Value = $u.Value }

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    @{Cell="H16"; Value=2112.2}
    @{Cell="I16"; Value=1861.909}
    @{Cell="J16"; Value=2800.5}
    @{Cell="K16"; Value=1861.909}
    @{Cell="L16"; Value=2800.5}
    @{Cell="M16"; Value=-1574.909}
    @{Cell="N16"; Value=-3374.5}
    @{Cell="H19"; Value=339.6}
    @{Cell="I19"; Value=339.6}
    @{Cell="J19"; Value=0}
    @{Cell="K19"; Value=339.6}
    @{Cell="L19"; Value=0}
    @{Cell="M19"; Value=-169.6}
    @{Cell="H24"; Value=339.6}
    @{Cell="I24"; Value=339.6}
    @{Cell="J24"; Value=0}
    @{Cell="K24"; Value=339.6}
    @{Cell="L24"; Value=0}
    @{Cell="M24"; Value=-169.6}
    @{Cell="H31"; Value=630930.1}
    @{Cell="I31"; Value=4746.654}
    @{Cell="J31"; Value=1124286.8}
    @{Cell="K31"; Value=4746.654}
    @{Cell="L31"; Value=1124286.8}
    @{Cell="M31"; Value=-4451.654}
    @{Cell="N31"; Value=-1124876.8}
    @{Cell="H34"; Value=630930.1}
    @{Cell="I34"; Value=4746.654}
    @{Cell="J34"; Value=1124286.8}
    @{Cell="K34"; Value=4746.654}
    @{Cell="L34"; Value=1124286.8}
    @{Cell="M34"; Value=-4544.654}
    @{Cell="N34"; Value=-1124690.8}
    @{Cell="H43"; Value=21472}
    @{Cell="J43"; Value=21472}
    @{Cell="L43"; Value=21472}
    @{Cell="N43"; Value=-21840}
    @{Cell="H58"; Value=4792688}
    @{Cell="I58"; Value=6995280}
    @{Cell="J58"; Value=20405.666}
    @{Cell="K58"; Value=6995280}
    @{Cell="L58"; Value=20405.666}
    @{Cell="M58"; Value=-6995077}
    @{Cell="N58"; Value=-20811.666}
    @{Cell="H101"; Value=21472}
    @{Cell="J101"; Value=21472}
    @{Cell="L101"; Value=21472}
    @{Cell="N101"; Value=-27962}
    @{Cell="H107"; Value=749.1818}
    @{Cell="I107"; Value=737.2381}
    @{Cell="K107"; Value=737.2381}
    @{Cell="M107"; Value=1182.7619}
    @{Cell="H113"; Value=2112.2}
    @{Cell="I113"; Value=1861.909}
    @{Cell="J113"; Value=2800.5}
    @{Cell="K113"; Value=1861.909}
    @{Cell="L113"; Value=2800.5}
    @{Cell="M113"; Value=308.0909999999999}
    @{Cell="N113"; Value=-7140.5}
    @{Cell="H136"; Value=4792688}
    @{Cell="I136"; Value=6995280}
    @{Cell="J136"; Value=20405.666}
    @{Cell="K136"; Value=20985840}
    @{Cell="L136"; Value=61216.99800000001}
    @{Cell="M136"; Value=-20983290}
    @{Cell="N136"; Value=-66316.99800000001}
)
foreach ($u in $updates) { $ws.Range($u.Cell).Value = $u.Value }
foreach ($c in @("N19", "N24")) { $ws.Range($c).ClearContents() }

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    @{Cell="H25"; Value=2475.25}
    @{Cell="I25"; Value=1950.5}
    @{Cell="J25"; Value=3000}
    @{Cell="K25"; Value=5851.5}
    @{Cell="L25"; Value=9000}
    @{Cell="M25"; Value=-5682.5}
    @{Cell="N25"; Value=-9338}
    @{Cell="H30"; Value=2475.25}
    @{Cell="I30"; Value=1950.5}
    @{Cell="J30"; Value=3000}
    @{Cell="K30"; Value=5851.5}
    @{Cell="L30"; Value=9000}
    @{Cell="M30"; Value=-5749.5}
    @{Cell="N30"; Value=-9204}
    @{Cell="H33"; Value=109}
    @{Cell="I33"; Value=85}
    @{Cell="J33"; Value=125}
    @{Cell="K33"; Value=510}
    @{Cell="L33"; Value=750}
    @{Cell="M33"; Value=-227}
    @{Cell="N33"; Value=-1316}
    @{Cell="H107"; Value=627.12964}
    @{Cell="J107"; Value=1862.5}
    @{Cell="L107"; Value=5587.5}
    @{Cell="N107"; Value=-9427.5}
    @{Cell="H122"; Value=798.125}
    @{Cell="I122"; Value=219}
    @{Cell="J122"; Value=1542.7142}
    @{Cell="K122"; Value=1971}
    @{Cell="L122"; Value=13884.4278}
    @{Cell="M122"; Value=479}
    @{Cell="N122"; Value=-18784.4278}
)
foreach ($u in $updates) { $ws.Range($u.Cell).Value = $u.Value }

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    @{Cell="H41"; Value=4234.6}
    @{Cell="I41"; Value=1620.8572}
    @{Cell="J41"; Value=10333.333}
    @{Cell="K41"; Value=1620.8572}
    @{Cell="L41"; Value=10333.333}
    @{Cell="M41"; Value=-1265.8572}
    @{Cell="N41"; Value=-11043.333}
)
foreach ($u in $updates) { $ws.Range($u.Cell).Value = $u.Value }

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    @{Cell="H140"; Value=60465.727}
    @{Cell="J140"; Value=60465.727}
    @{Cell="L140"; Value=60465.727}
    @{Cell="N140"; Value=-70825.727}
)
foreach ($u in $updates) { $ws.Range($u.Cell).Value = $u.Value }

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    @{Cell="H107"; Value=1647.1765}
    @{Cell="I107"; Value=1096.125}
    @{Cell="J107"; Value=2137}
    @{Cell="K107"; Value=3288.375}
    @{Cell="L107"; Value=6411}
    @{Cell="M107"; Value=-1368.375}
    @{Cell="N107"; Value=-10251}
)
foreach ($u in $updates) { $ws.Range($u.Cell).Value = $u.Value }
